$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for every existing data row (2..397) from 45192 -> 45202
for ($r = 2; $r -le 397; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45202
}

# Append new row 398 with the new case data
$ws.Cells.Item(398, 1).Value2 = "A 45586-2023"
$ws.Cells.Item(398, 2).Value2 = 45194
$ws.Cells.Item(398, 3).Value2 = 45202
$ws.Cells.Item(398, 4).Value2 = "HALLANDS LÄN"
$ws.Cells.Item(398, 5).Value2 = "HALMSTAD"
$ws.Cells.Item(398, 7).Value2 = 1.5
$ws.Cells.Item(398, 8).Value2 = 0
$ws.Cells.Item(398, 9).Value2 = 0
$ws.Cells.Item(398, 10).Value2 = 0
$ws.Cells.Item(398, 11).Value2 = 0
$ws.Cells.Item(398, 12).Value2 = 0
$ws.Cells.Item(398, 13).Value2 = 0
$ws.Cells.Item(398, 14).Value2 = 0
$ws.Cells.Item(398, 15).Value2 = 0
$ws.Cells.Item(398, 16).Value2 = 0
$ws.Cells.Item(398, 17).Value2 = 0

# Match styles used by the rest of the table
$ws.Cells.Item(398, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(398, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(398, 18).WrapText = $true

$ws.Rows.Item(397).RowHeight = 15
$ws.Rows.Item(398).RowHeight = 15
